# Auto-generated Excel COM-interop script applying the Lamia_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2515.7058
$ws.Range("I33").Value = 104.92857
$ws.Range("K33").Value = 104.92857
$ws.Range("M33").Value = 124.07143
# Row 64
$ws.Range("H64").Value = 9584.200000000001
$ws.Range("I64").Value = 4499
$ws.Range("J64").Value = 9796.083000000001
$ws.Range("K64").Value = 4499
$ws.Range("L64").Value = 9796.083000000001
$ws.Range("M64").Value = -4251
$ws.Range("N64").Value = -10292.083
# Row 67
$ws.Range("H67").Value = 9584.200000000001
$ws.Range("I67").Value = 4499
$ws.Range("J67").Value = 9796.083000000001
$ws.Range("K67").Value = 4499
$ws.Range("L67").Value = 9796.083000000001
$ws.Range("M67").Value = -3641
$ws.Range("N67").Value = -11512.083
# Row 132
$ws.Range("H132").Value = 1918.2391
$ws.Range("I132").Value = 1999.878
$ws.Range("K132").Value = 5999.634
$ws.Range("M132").Value = -3469.634
# Row 133
$ws.Range("H133").Value = 69999
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("N133").Value = -80119
# Row 135
$ws.Range("H135").Value = 604.4666999999999
$ws.Range("I135").Value = 650.46155
$ws.Range("K135").Value = 5854.15395
$ws.Range("M135").Value = -3319.15395
# Row 136
$ws.Range("H136").Value = 68993.60000000001
$ws.Range("J136").Value = 68993.60000000001
$ws.Range("L136").Value = 68993.60000000001
$ws.Range("N136").Value = -79193.60000000001
# Row 138
$ws.Range("H138").Value = 5512.96
$ws.Range("I138").Value = 4169.5713
$ws.Range("K138").Value = 12508.7139
$ws.Range("M138").Value = -7368.713899999999
# Row 139
$ws.Range("H139").Value = 59997.5
$ws.Range("J139").Value = 69996.664
$ws.Range("L139").Value = 69996.664
$ws.Range("N139").Value = -80276.664
# Row 140
$ws.Range("H140").Value = 67715.39999999999
$ws.Range("J140").Value = 70239.336
$ws.Range("L140").Value = 70239.336
$ws.Range("N140").Value = -80599.336

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 8312.833000000001
$ws.Range("I2").Value = 3849
$ws.Range("J2").Value = 9588.214
$ws.Range("K2").Value = 3849
$ws.Range("L2").Value = 9588.214
$ws.Range("M2").Value = -3736
$ws.Range("N2").Value = -9814.214
# Row 4
$ws.Range("H4").Value = 2067.875
$ws.Range("I4").Value = 181.66667
$ws.Range("K4").Value = 181.66667
$ws.Range("M4").Value = -65.66667000000001
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 102
$ws.Range("H102").Value = 1447.2
$ws.Range("I102").Value = 1447.2
$ws.Range("K102").Value = 1447.2
$ws.Range("M102").Value = 174.8
# Row 116
$ws.Range("H116").Value = 8312.833000000001
$ws.Range("I116").Value = 3849
$ws.Range("J116").Value = 9588.214
$ws.Range("K116").Value = 3849
$ws.Range("L116").Value = 9588.214
$ws.Range("M116").Value = -1555
$ws.Range("N116").Value = -14176.214
# Row 132
$ws.Range("H132").Value = 2372.5881
$ws.Range("J132").Value = 4471.1665
$ws.Range("L132").Value = 13413.4995
$ws.Range("N132").Value = -18473.4995

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 8312.833000000001
$ws.Range("I3").Value = 3849
$ws.Range("J3").Value = 9588.214
$ws.Range("K3").Value = 3849
$ws.Range("L3").Value = 9588.214
$ws.Range("M3").Value = -3735
$ws.Range("N3").Value = -9816.214
# Row 107
$ws.Range("H107").Value = 1929.8572
$ws.Range("I107").Value = 1252.3334
$ws.Range("K107").Value = 1252.3334
$ws.Range("M107").Value = 667.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 287.85715
$ws.Range("I7").Value = 106
$ws.Range("J7").Value = 388.8889
$ws.Range("K7").Value = 106
$ws.Range("L7").Value = 388.8889
$ws.Range("M7").Value = 7
$ws.Range("N7").Value = -614.8888999999999
# Row 105
$ws.Range("H105").Value = 6037.037
$ws.Range("I105").Value = 6799.364
$ws.Range("J105").Value = 5512.9375
$ws.Range("K105").Value = 6799.364
$ws.Range("L105").Value = 5512.9375
$ws.Range("M105").Value = -5052.364
$ws.Range("N105").Value = -9006.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1438.1666
$ws.Range("I107").Value = 1406.4286
$ws.Range("J107").Value = 1549.25
$ws.Range("K107").Value = 4219.2858
$ws.Range("L107").Value = 4647.75
$ws.Range("M107").Value = -2299.2858
$ws.Range("N107").Value = -8487.75
# Row 114
$ws.Range("H114").Value = 125002130
$ws.Range("I114").Value = 2263.5
$ws.Range("J114").Value = 250002000
$ws.Range("K114").Value = 6790.5
$ws.Range("L114").Value = 750006000
$ws.Range("M114").Value = -3536.5
$ws.Range("N114").Value = -750012508
# Row 132
$ws.Range("H132").Value = 5583.5557
$ws.Range("I132").Value = 4400.8
$ws.Range("K132").Value = 39607.2
$ws.Range("M132").Value = -37077.2

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 4177.28
$ws.Range("I113").Value = 1582.1428
$ws.Range("K113").Value = 1582.1428
$ws.Range("M113").Value = 587.8571999999999
# Row 122
$ws.Range("H122").Value = 8803.625
$ws.Range("I122").Value = 7962.6313
$ws.Range("K122").Value = 23887.8939
$ws.Range("M122").Value = -21437.8939
# Row 130
$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040
# Row 132
$ws.Range("H132").Value = 4457.143
$ws.Range("I132").Value = 1960.3334
$ws.Range("K132").Value = 5881.0002
$ws.Range("M132").Value = -3351.0002

$ws = $wb.Worksheets.Item("LTW")
# Row 106
$ws.Range("H106").Value = 8685
$ws.Range("J106").Value = 8685
$ws.Range("L106").Value = 8685
$ws.Range("N106").Value = -11209
# Row 122
$ws.Range("H122").Value = 6662.3125
$ws.Range("I122").Value = 4186.5454
$ws.Range("K122").Value = 12559.6362
$ws.Range("M122").Value = -10109.6362
# Row 136
$ws.Range("H136").Value = 7863.875
$ws.Range("I136").Value = 5742.92
$ws.Range("K136").Value = 17228.76
$ws.Range("M136").Value = -14678.76

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2908.1667
$ws.Range("J96").Value = 2333.3333
$ws.Range("L96").Value = 2333.3333
$ws.Range("N96").Value = -5079.3333
# Row 100
$ws.Range("H100").Value = 683.82355
$ws.Range("I100").Value = 705.7857
$ws.Range("J100").Value = 581.3333
$ws.Range("K100").Value = 1411.5714
$ws.Range("L100").Value = 1162.6666
$ws.Range("M100").Value = -870.5714
$ws.Range("N100").Value = -2244.6666
# Row 113
$ws.Range("H113").Value = 1457.8889
$ws.Range("I113").Value = 1668.091
$ws.Range("J113").Value = 1127.5714
$ws.Range("K113").Value = 5004.272999999999
$ws.Range("L113").Value = 3382.7142
$ws.Range("M113").Value = -2834.272999999999
$ws.Range("N113").Value = -7722.7142
# Row 122
$ws.Range("H122").Value = 4818.0415
$ws.Range("I122").Value = 3035
$ws.Range("K122").Value = 9105
$ws.Range("M122").Value = -6655
# Row 125
$ws.Range("H125").Value = 88325.86
$ws.Range("J125").Value = 99713.5
$ws.Range("L125").Value = 99713.5
$ws.Range("N125").Value = -109553.5
# Row 132
$ws.Range("H132").Value = 3879.9048
$ws.Range("I132").Value = 2885.4412
$ws.Range("K132").Value = 8656.3236
$ws.Range("M132").Value = -6126.3236
# Row 136
$ws.Range("H136").Value = 3352.8857
$ws.Range("I136").Value = 1862.5714
$ws.Range("K136").Value = 5587.7142
$ws.Range("M136").Value = -3037.7142
